$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" target cluster no longer exists in the refreshed TPM run,
# so the 3 rows that referenced it (rows 6, 11, 16 in the old layout) are gone and
# every remaining sending/target-cluster pair is rewritten with the new values.
# Drop the 3 now-obsolete trailing rows so the table is 3 senders x 4 targets = 12 rows.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lgi2"
$ws.Cells.Item(2, 3).Value = "Adam11"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.1376636666666667
$ws.Cells.Item(2, 8).Value = 0.412991
$ws.Cells.Item(2, 9).Value = 0.01821680097623009
$ws.Cells.Item(2, 10).Value = 0.01821680097623009
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1266143333333333
$ws.Cells.Item(2, 14).Value = 0.379843
$ws.Cells.Item(2, 15).Value = 0.0955140869844225
$ws.Cells.Item(2, 16).Value = 0.0955140869844225
$ws.Cells.Item(2, 17).Value = 0.01743019337922222
$ws.Cells.Item(2, 18).Value = 0.156871740413
$ws.Cells.Item(2, 19).Value = 0.001739961113021553
$ws.Cells.Item(2, 20).Value = 0.001739961113021553

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lgi2"
$ws.Cells.Item(3, 3).Value = "Adam11"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.1376636666666667
$ws.Cells.Item(3, 8).Value = 0.412991
$ws.Cells.Item(3, 9).Value = 0.01821680097623009
$ws.Cells.Item(3, 10).Value = 0.01821680097623009
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3243313333333334
$ws.Cells.Item(3, 14).Value = 0.972994
$ws.Cells.Item(3, 15).Value = 0.2446659107876707
$ws.Cells.Item(3, 16).Value = 0.2446659107876707
$ws.Cells.Item(3, 17).Value = 0.04464864056155556
$ws.Cells.Item(3, 18).Value = 0.401837765054
$ws.Cells.Item(3, 19).Value = 0.004457030202487063
$ws.Cells.Item(3, 20).Value = 0.004457030202487063

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lgi2"
$ws.Cells.Item(4, 3).Value = "Adam11"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.1376636666666667
$ws.Cells.Item(4, 8).Value = 0.412991
$ws.Cells.Item(4, 9).Value = 0.01821680097623009
$ws.Cells.Item(4, 10).Value = 0.01821680097623009
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.05644366666666667
$ws.Cells.Item(4, 14).Value = 0.169331
$ws.Cells.Item(4, 15).Value = 0.0425794232437066
$ws.Cells.Item(4, 16).Value = 0.0425794232437066
$ws.Cells.Item(4, 17).Value = 0.007770242113444444
$ws.Cells.Item(4, 18).Value = 0.069932179021
$ws.Cells.Item(4, 19).Value = 0.0007756608789132685
$ws.Cells.Item(4, 20).Value = 0.0007756608789132685

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Lgi2"
$ws.Cells.Item(5, 3).Value = "Adam11"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1376636666666667
$ws.Cells.Item(5, 8).Value = 0.412991
$ws.Cells.Item(5, 9).Value = 0.01821680097623009
$ws.Cells.Item(5, 10).Value = 0.01821680097623009
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8182196666666667
$ws.Cells.Item(5, 14).Value = 2.454659
$ws.Cells.Item(5, 15).Value = 0.6172405789842003
$ws.Cells.Item(5, 16).Value = 0.6172405789842002
$ws.Cells.Item(5, 17).Value = 0.1126391194521111
$ws.Cells.Item(5, 18).Value = 1.013752075069
$ws.Cells.Item(5, 19).Value = 0.0112441487818082
$ws.Cells.Item(5, 20).Value = 0.0112441487818082

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lgi2"
$ws.Cells.Item(6, 3).Value = "Adam11"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.314644999999999
$ws.Cells.Item(6, 8).Value = 21.943935
$ws.Cells.Item(6, 9).Value = 0.9679346439276632
$ws.Cells.Item(6, 10).Value = 0.967934643927663
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1266143333333333
$ws.Cells.Item(6, 14).Value = 0.379843
$ws.Cells.Item(6, 15).Value = 0.0955140869844225
$ws.Cells.Item(6, 16).Value = 0.0955140869844225
$ws.Cells.Item(6, 17).Value = 0.9261389002449998
$ws.Cells.Item(6, 18).Value = 8.335250102204999
$ws.Cells.Item(6, 19).Value = 0.09245139377534284
$ws.Cells.Item(6, 20).Value = 0.09245139377534282

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lgi2"
$ws.Cells.Item(7, 3).Value = "Adam11"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.314644999999999
$ws.Cells.Item(7, 8).Value = 21.943935
$ws.Cells.Item(7, 9).Value = 0.9679346439276632
$ws.Cells.Item(7, 10).Value = 0.967934643927663
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3243313333333334
$ws.Cells.Item(7, 14).Value = 0.972994
$ws.Cells.Item(7, 15).Value = 0.2446659107876707
$ws.Cells.Item(7, 16).Value = 0.2446659107876707
$ws.Cells.Item(7, 17).Value = 2.37236856571
$ws.Cells.Item(7, 18).Value = 21.35131709139
$ws.Cells.Item(7, 19).Value = 0.2368206112395014
$ws.Cells.Item(7, 20).Value = 0.2368206112395014

# Row 8: FAPs -> Inflammatory-Mac
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Lgi2"
$ws.Cells.Item(8, 3).Value = "Adam11"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.314644999999999
$ws.Cells.Item(8, 8).Value = 21.943935
$ws.Cells.Item(8, 9).Value = 0.9679346439276632
$ws.Cells.Item(8, 10).Value = 0.967934643927663
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.05644366666666667
$ws.Cells.Item(8, 14).Value = 0.169331
$ws.Cells.Item(8, 15).Value = 0.0425794232437066
$ws.Cells.Item(8, 16).Value = 0.0425794232437066
$ws.Cells.Item(8, 17).Value = 0.412865384165
$ws.Cells.Item(8, 18).Value = 3.715788457485
$ws.Cells.Item(8, 19).Value = 0.04121409887604242
$ws.Cells.Item(8, 20).Value = 0.04121409887604241

# Row 9: FAPs -> MuSCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Lgi2"
$ws.Cells.Item(9, 3).Value = "Adam11"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.314644999999999
$ws.Cells.Item(9, 8).Value = 21.943935
$ws.Cells.Item(9, 9).Value = 0.9679346439276632
$ws.Cells.Item(9, 10).Value = 0.967934643927663
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.8182196666666667
$ws.Cells.Item(9, 14).Value = 2.454659
$ws.Cells.Item(9, 15).Value = 0.6172405789842003
$ws.Cells.Item(9, 16).Value = 0.6172405789842002
$ws.Cells.Item(9, 17).Value = 5.984986393684999
$ws.Cells.Item(9, 18).Value = 53.86487754316499
$ws.Cells.Item(9, 19).Value = 0.5974485400367765
$ws.Cells.Item(9, 20).Value = 0.5974485400367764

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Lgi2"
$ws.Cells.Item(10, 3).Value = "Adam11"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.104653
$ws.Cells.Item(10, 8).Value = 0.313959
$ws.Cells.Item(10, 9).Value = 0.01384855509610675
$ws.Cells.Item(10, 10).Value = 0.01384855509610675
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1266143333333333
$ws.Cells.Item(10, 14).Value = 0.379843
$ws.Cells.Item(10, 15).Value = 0.0955140869844225
$ws.Cells.Item(10, 16).Value = 0.0955140869844225
$ws.Cells.Item(10, 17).Value = 0.01325056982633333
$ws.Cells.Item(10, 18).Value = 0.119255128437
$ws.Cells.Item(10, 19).Value = 0.001322732096058108
$ws.Cells.Item(10, 20).Value = 0.001322732096058108

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Lgi2"
$ws.Cells.Item(11, 3).Value = "Adam11"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.104653
$ws.Cells.Item(11, 8).Value = 0.313959
$ws.Cells.Item(11, 9).Value = 0.01384855509610675
$ws.Cells.Item(11, 10).Value = 0.01384855509610675
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.3243313333333334
$ws.Cells.Item(11, 14).Value = 0.972994
$ws.Cells.Item(11, 15).Value = 0.2446659107876707
$ws.Cells.Item(11, 16).Value = 0.2446659107876707
$ws.Cells.Item(11, 17).Value = 0.03394224702733333
$ws.Cells.Item(11, 18).Value = 0.305480223246
$ws.Cells.Item(11, 19).Value = 0.003388269345682196
$ws.Cells.Item(11, 20).Value = 0.003388269345682196

# Row 12: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Lgi2"
$ws.Cells.Item(12, 3).Value = "Adam11"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.104653
$ws.Cells.Item(12, 8).Value = 0.313959
$ws.Cells.Item(12, 9).Value = 0.01384855509610675
$ws.Cells.Item(12, 10).Value = 0.01384855509610675
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.05644366666666667
$ws.Cells.Item(12, 14).Value = 0.169331
$ws.Cells.Item(12, 15).Value = 0.0425794232437066
$ws.Cells.Item(12, 16).Value = 0.0425794232437066
$ws.Cells.Item(12, 17).Value = 0.005906999047666667
$ws.Cells.Item(12, 18).Value = 0.053162991429
$ws.Cells.Item(12, 19).Value = 0.0005896634887509192
$ws.Cells.Item(12, 20).Value = 0.0005896634887509192

# Row 13: MuSCs -> MuSCs
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Lgi2"
$ws.Cells.Item(13, 3).Value = "Adam11"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.104653
$ws.Cells.Item(13, 8).Value = 0.313959
$ws.Cells.Item(13, 9).Value = 0.01384855509610675
$ws.Cells.Item(13, 10).Value = 0.01384855509610675
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.8182196666666667
$ws.Cells.Item(13, 14).Value = 2.454659
$ws.Cells.Item(13, 15).Value = 0.6172405789842003
$ws.Cells.Item(13, 16).Value = 0.6172405789842002
$ws.Cells.Item(13, 17).Value = 0.08562914277566666
$ws.Cells.Item(13, 18).Value = 0.770662284981
$ws.Cells.Item(13, 19).Value = 0.008547890165615527
$ws.Cells.Item(13, 20).Value = 0.008547890165615524

